$wb = $excel.ActiveWorkbook

# --- Brand sheet (sheet1.xml) ---
# Add two new rows with "Name" values completing the product/brand combobox list.
$wsBrand = $wb.Worksheets.Item("Brand")
$wsBrand.Range("A4").Value = 3
$wsBrand.Range("B4").Value = "dfbdfb"
$wsBrand.Range("A5").Value = 4
# "1" must stay a text value (not be coerced to a number), so format the cell
# as Text before typing the numeric-looking label.
$wsBrand.Range("B5").NumberFormat = "@"
$wsBrand.Range("B5").Value = "1"

# --- Client sheet (sheet2.xml) ---
# Add a new row with repeated "11" values across columns B, C and D.
$wsClient = $wb.Worksheets.Item("Client")
$wsClient.Range("A3").Value = 2
$wsClient.Range("B3:D3").NumberFormat = "@"
$wsClient.Range("B3").Value = "11"
$wsClient.Range("C3").Value = "11"
$wsClient.Range("D3").Value = "11"
